# Update column F (dSF) values to reflect repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -2
    4  = -2
    5  = 2
    6  = -2
    7  = 3
    8  = -3
    9  = 1
    10 = -1
    13 = 3
    14 = 1
    15 = -6
    16 = 4
    17 = -2
    18 = -2
    19 = -3
    20 = 1
    21 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
